# Weekly update: insert 2 new price rows for "Carson" variety (Especial/Primera)
# at the top of the Durazno block (rows 124-125), pushing all subsequent
# rows down by 2 (124-191 -> 126-193). New rows carry updated date,
# volume, price-range and unit-size figures for the latest week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 124.. down by two to make room for the new entries.
$ws.Rows.Item(124).Resize(2).Insert()

$ws.Cells.Item(124, 1).Value = 4
$ws.Cells.Item(124, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(124, 3).Value = 'Los Lagos'
$ws.Cells.Item(124, 4).Value = 44603
$ws.Cells.Item(124, 5).Value = 10
$ws.Cells.Item(124, 6).Value = 'Fruta'
$ws.Cells.Item(124, 7).Value = 100103
$ws.Cells.Item(124, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(124, 9).Value = 100103004
$ws.Cells.Item(124, 10).Value = 'Durazno'
$ws.Cells.Item(124, 11).Value = 'Carson'
$ws.Cells.Item(124, 12).Value = 'Especial'
$ws.Cells.Item(124, 13).Value = 200
$ws.Cells.Item(124, 14).Value = 20000
$ws.Cells.Item(124, 15).Value = 20000
$ws.Cells.Item(124, 16).Value = 20000
$ws.Cells.Item(124, 17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(124, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(124, 19).Value = 1333
$ws.Cells.Item(124, 20).Value = 15
$ws.Cells.Item(125, 1).Value = 4
$ws.Cells.Item(125, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(125, 3).Value = 'Los Lagos'
$ws.Cells.Item(125, 4).Value = 44603
$ws.Cells.Item(125, 5).Value = 10
$ws.Cells.Item(125, 6).Value = 'Fruta'
$ws.Cells.Item(125, 7).Value = 100103
$ws.Cells.Item(125, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(125, 9).Value = 100103004
$ws.Cells.Item(125, 10).Value = 'Durazno'
$ws.Cells.Item(125, 11).Value = 'Carson'
$ws.Cells.Item(125, 12).Value = 'Primera'
$ws.Cells.Item(125, 13).Value = 400
$ws.Cells.Item(125, 14).Value = 16000
$ws.Cells.Item(125, 15).Value = 17000
$ws.Cells.Item(125, 16).Value = 16500
$ws.Cells.Item(125, 17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(125, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(125, 19).Value = 1100
$ws.Cells.Item(125, 20).Value = 15

Write-Host "Inserted 2 rows and populated new row 124/125 data"
